$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.091562213532490944, 0.091028556423268014),
    @(-0.084611713586231652, 0.083058099347446124),
    @(-0.037443331664029245, 0.037017796224608546),
    @(-0.029017796274796837, 0.028634991253118969),
    @(-0.025634991276170638, 0.024338301957061326),
    @(0.013789195129410459, -0.014157178039539531),
    @(0.024157177976674049, -0.024250456186025282),
    @(0.034250456124968576, -0.034448574663360709),
    @(0.036448574651156473, -0.036623952737774967),
    @(-0.022894974363216747, 0.022891325900738835),
    @(-0.019891325913725666, 0.019885276398429852),
    @(-0.016385276414815575, 0.016348072238788536),
    @(-0.012848072257239274, 0.012837086091498584),
    @(-0.0048370861355362393, 0.0048363255534500382),
    @(-0.008052971055660052, 0.0080344731301833505),
    @(-0.0060344731412782515, 0.0060034335598340149),
    @(-0.0040034335715102287, 0.0039999999769202432),
    @(-0.00041684287132071063, 0.0003114689895866718),
    @(0.0036885309873251337, -0.0044736115525885189),
    @(0.0084736115298955639, -0.0087015745101517439),
    @(-0.0040057101807455808, 0.0039999999773687733),
    @(-0.041614767710179024, 0.041362477994618452),
    @(-0.04050270317412874, 0.040099743070967975),
    @(-0.020099743192176334, 0.019999999877073904),
    @(-0.036127497105326256, 0.036091664818602709),
    @(-0.033591664840830759, 0.033548996006988574),
    @(-0.03104899602999911, 0.030816616075918901),
    @(-0.028816616099154757, 0.028676147931462559),
    @(-0.021676147984823757, 0.02164806514892792),
    @(0.038351934500494345, -0.038461712356157474),
    @(0.045461712305778335, -0.04551151334352177),
    @(0.055511513277785696, -0.055639658857961294)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
